$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.082.49"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.836.99"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'243.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").Value = "'0.6288"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.32%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.07567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.91%  "

$ws.Range("D9").Value = "'0.2937"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").Value = "'22.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").Value = "'0.07753"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "1.840.54"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").Value = "'4.974"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "'0.6667"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "'0.00001003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +15.48%  "

$ws.Range("D16").Value = "'83.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.21%  "

$ws.Range("D17").Value = "'6.083"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").Value = "29.101.71"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "'226.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").Value = "'12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").Value = "'159.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.94%  "

$ws.Range("E25").Value = "  +1.33%  "

$ws.Range("D26").Value = "'8.508"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "'17.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").Value = "'1.498"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").Value = "'4.102"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "'4.014"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.54%  "

$ws.Range("D31").Value = "'1.195"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("D32").Value = "'0.05263"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("D33").Value = "'1.850"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("D34").Value = "'0.7379"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").Value = "'1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").Value = "'2.681"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("D37").Value = "1.246.34"
$ws.Range("E37").Value = "  -4.33%  "

$ws.Range("D38").Value = "'2.762"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("D39").Value = "'0.01788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").Value = "'6.382"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("D41").Value = "'0.9024"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "'1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").Value = "'102.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").Value = "'0.00000000128"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.39%  "

$ws.Range("D45").Value = "1.988.99"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("D46").Value = "'64.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").Value = "'0.5122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").Value = "'0.4047"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").Value = "'8.904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.43%  "

$ws.Range("D50").Value = "'0.05771"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").Value = "'6.722"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
